$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "B2" "2.0"
Set-TextValue "C2" "-9.85203026391962"
Set-TextValue "D2" "1.04848751138525"

# Row 3
Set-TextValue "B3" "0.951512488614748"
Set-TextValue "C3" "-5.70057024136476"
Set-TextValue "D3" "0.428928774956478"

# Row 4
Set-TextValue "B4" "0.52258371365827"
Set-TextValue "C4" "-1.93462565319986"
Set-TextValue "D4" "0.117995951340695"

# Row 5
Set-TextValue "B5" "0.404587762317575"
Set-TextValue "C5" "-0.261206030131537"
Set-TextValue "D5" "0.0153770222176567"

# Row 6
Set-TextValue "B6" "0.389210740099918"
Set-TextValue "C6" "-0.0053489238075823"
Set-TextValue "D6" "0.0003146487342669"

# New row 7
Set-TextValue "A7" "6"
Set-TextValue "B7" "0.388896091365651"
Set-TextValue "C7" "-2.2901295766697e-06"
Set-TextValue "D7" "1.34716096833465e-07"
